$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (image) values for rows 2..33, in row order
$bValues = @(
    "flower/flower093.png",
    "flower/flower097.png",
    "flower/flower082.png",
    "face/face119.png",
    "face/face096.png",
    "face/face091.png",
    "flower/flower109.png",
    "face/face102.png",
    "face/face092.png",
    "face/face104.png",
    "face/face071.png",
    "face/face098.png",
    "flower/flower092.png",
    "face/face069.png",
    "flower/flower091.png",
    "flower/flower094.png",
    "flower/flower087.png",
    "face/face090.png",
    "face/face083.png",
    "face/face097.png",
    "flower/flower068.png",
    "flower/flower122.png",
    "face/face080.png",
    "flower/flower083.png",
    "face/face115.png",
    "flower/flower069.png",
    "flower/flower072.png",
    "face/face114.png",
    "flower/flower070.png",
    "face/face094.png",
    "flower/flower065.png",
    "flower/flower114.png"
)

# Column C (word) values for rows 2..33, in row order
$cValues = @(
    "hauen",
    "drehen",
    "töten",
    "pflegen",
    "klappen",
    "husten",
    "dauern",
    "starten",
    "antun",
    "biegen",
    "lehnen",
    "schicken",
    "rücken",
    "loben",
    "stechen",
    "scheitern",
    "rasen",
    "schenken",
    "hupen",
    "fliegen",
    "regnen",
    "fesseln",
    "segeln",
    "mieten",
    "tauschen",
    "langen",
    "sondern",
    "krachen",
    "gründen",
    "stärken",
    "füllen",
    "fühlen"
)

# Column D (category) values for rows 2..33, in row order
$dValues = @(
    "flower",
    "flower",
    "flower",
    "face",
    "face",
    "face",
    "flower",
    "face",
    "face",
    "face",
    "face",
    "face",
    "flower",
    "face",
    "flower",
    "flower",
    "flower",
    "face",
    "face",
    "face",
    "flower",
    "flower",
    "face",
    "flower",
    "face",
    "flower",
    "flower",
    "face",
    "flower",
    "face",
    "flower",
    "flower"
)

$startRow = 2

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $bValues[$i]
}

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $cValues[$i]
}

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $dValues[$i]
}

